$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Masterdata refresh (2nd May) - update machine id values
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Selection moved to full rows starting at 34 (e.g. user selected row 34 downward)
$ws.Rows("34:1048576").Select()
